$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Акция" (promotion / sale flag) column -------------------------
# Header: reuse the same look as the other header cells (bold, centered).
$ws.Range("D1").Copy() | Out-Null
$ws.Range("E1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("E1").Value = "Акция"

# Body cells: default font, centered. Wines currently on sale get "+",
# the rest stay blank (still formatted/centered).
$body = $ws.Range("E2:E7")
$body.Font.Name = "Arial"
$body.Font.Size = 10
$body.Font.Bold = $false
$body.HorizontalAlignment = -4108   # xlCenter

$ws.Range("E2").Value = ""
$ws.Range("E4").Value = ""
$ws.Range("E5").Value = ""
$ws.Range("E7").Value = ""

# "Гранатовый браслет" and "Ркацители" are on sale -> mark with "+".
# These use the regular body font (11pt, matches column A-D data rows).
$ws.Range("E3").Font.Name = "Arial"
$ws.Range("E3").Font.Size = 11
$ws.Range("E3").Font.Bold = $false
$ws.Range("E3").HorizontalAlignment = -4108
$ws.Range("E3").Value = "+"

$ws.Range("E6").Font.Name = "Arial"
$ws.Range("E6").Font.Size = 11
$ws.Range("E6").Font.Bold = $false
$ws.Range("E6").HorizontalAlignment = -4108
$ws.Range("E6").Value = "+"
